$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 4428.5
$ws.Range("I2").Value = 5194
$ws.Range("K2").Value = 5194
$ws.Range("M2").Value = -5081
$ws.Range("H7").Value = 7084.1665
$ws.Range("I7").Value = 5
$ws.Range("K7").Value = 5
$ws.Range("M7").Value = 107
$ws.Range("H14").Value = 7084.1665
$ws.Range("I14").Value = 5
$ws.Range("K14").Value = 5
$ws.Range("M14").Value = 186
$ws.Range("H17").Value = 1066.3
$ws.Range("J17").Value = 1066.3
$ws.Range("L17").Value = 3198.9
$ws.Range("N17").Value = -3534.9
$ws.Range("H58").Value = 785125.4
$ws.Range("I58").Value = 980656.75
$ws.Range("J58").Value = 3000
$ws.Range("K58").Value = 2941970.25
$ws.Range("L58").Value = 9000
$ws.Range("M58").Value = -2941820.25
$ws.Range("N58").Value = -9300
$ws.Range("H64").Value = 58313.055
$ws.Range("I64").Value = 93520
$ws.Range("J64").Value = 2987.8572
$ws.Range("K64").Value = 93520
$ws.Range("L64").Value = 2987.8572
$ws.Range("M64").Value = -93272
$ws.Range("N64").Value = -3483.8572
$ws.Range("H67").Value = 58313.055
$ws.Range("I67").Value = 93520
$ws.Range("J67").Value = 2987.8572
$ws.Range("K67").Value = 93520
$ws.Range("L67").Value = 2987.8572
$ws.Range("M67").Value = -92662
$ws.Range("N67").Value = -4703.8572
$ws.Range("H98").Value = 357.33334
$ws.Range("J98").Value = 499.66666
$ws.Range("L98").Value = 499.66666
$ws.Range("N98").Value = -3495.66666
$ws.Range("H122").Value = 357.33334
$ws.Range("J122").Value = 499.66666
$ws.Range("L122").Value = 1498.99998
$ws.Range("N122").Value = -6398.999980000001
$ws.Range("H127").Value = 38463070
$ws.Range("I127").Value = 484.14285
$ws.Range("J127").Value = 52633500
$ws.Range("K127").Value = 1452.42855
$ws.Range("L127").Value = 157900500
$ws.Range("M127").Value = 3507.57145
$ws.Range("N127").Value = -157910420
$ws.Range("H129").Value = 864.8
$ws.Range("J129").Value = 912.2273
$ws.Range("L129").Value = 2736.6819
$ws.Range("N129").Value = -12736.6819
$ws.Range("H137").Value = 1372.3226
$ws.Range("I137").Value = 1295.6296
$ws.Range("J137").Value = 1890
$ws.Range("K137").Value = 3886.8888
$ws.Range("L137").Value = 5670
$ws.Range("M137").Value = -1336.8888
$ws.Range("N137").Value = -10770

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H21").Value = 4763.4
$ws.Range("I21").Value = 950
$ws.Range("J21").Value = 20017
$ws.Range("K21").Value = 950
$ws.Range("L21").Value = 20017
$ws.Range("M21").Value = -576
$ws.Range("N21").Value = -20765
$ws.Range("H32").Value = 6755.923
$ws.Range("I32").Value = 5832.5967
$ws.Range("J32").Value = 25838
$ws.Range("K32").Value = 5832.5967
$ws.Range("L32").Value = 25838
$ws.Range("M32").Value = -5545.5967
$ws.Range("N32").Value = -26412
$ws.Range("H110").Value = 83426030
$ws.Range("I110").Value = 125138264
$ws.Range("J110").Value = 1575
$ws.Range("K110").Value = 125138264
$ws.Range("L110").Value = 1575
$ws.Range("M110").Value = -125136219
$ws.Range("N110").Value = -5665
$ws.Range("H122").Value = 1800.05
$ws.Range("I122").Value = 1744.3572
$ws.Range("J122").Value = 1930
$ws.Range("K122").Value = 5233.071599999999
$ws.Range("L122").Value = 5790
$ws.Range("M122").Value = -2783.071599999999
$ws.Range("N122").Value = -10690

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 612.7143
$ws.Range("I94").Value = 532.9231
$ws.Range("K94").Value = 532.9231
$ws.Range("M94").Value = -81.92309999999998
$ws.Range("H107").Value = 125001250
$ws.Range("I107").Value = 250000420
$ws.Range("K107").Value = 250000420
$ws.Range("M107").Value = -249998500

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H21").Value = 24800
$ws.Range("J21").Value = 24800
$ws.Range("L21").Value = 24800
$ws.Range("N21").Value = -25270
$ws.Range("H31").Value = 2270.4062
$ws.Range("I31").Value = 1653.9524
$ws.Range("J31").Value = 2571.465
$ws.Range("K31").Value = 1653.9524
$ws.Range("L31").Value = 2571.465
$ws.Range("M31").Value = -1358.9524
$ws.Range("N31").Value = -3161.465
$ws.Range("H34").Value = 2270.4062
$ws.Range("I34").Value = 1653.9524
$ws.Range("J34").Value = 2571.465
$ws.Range("K34").Value = 1653.9524
$ws.Range("L34").Value = 2571.465
$ws.Range("M34").Value = -1451.9524
$ws.Range("N34").Value = -2975.465
$ws.Range("H105").Value = 1136.3158
$ws.Range("I105").Value = 1078.4615
$ws.Range("J105").Value = 1261.6666
$ws.Range("K105").Value = 1078.4615
$ws.Range("L105").Value = 1261.6666
$ws.Range("M105").Value = 668.5385000000001
$ws.Range("N105").Value = -4755.6666

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1003.01
$ws.Range("J131").Value = 1025.7834
$ws.Range("L131").Value = 3077.3502
$ws.Range("N131").Value = -13157.3502

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 275587.72
$ws.Range("I102").Value = 2537.1
$ws.Range("K102").Value = 2537.1
$ws.Range("M102").Value = -915.0999999999999
$ws.Range("H132").Value = 2368.9644
$ws.Range("I132").Value = 2096.625
$ws.Range("J132").Value = 4003
$ws.Range("K132").Value = 6289.875
$ws.Range("L132").Value = 12009
$ws.Range("M132").Value = -3759.875
$ws.Range("N132").Value = -17069

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 85358.336
$ws.Range("I40").Value = 334100
$ws.Range("J40").Value = 2444.4443
$ws.Range("K40").Value = 334100
$ws.Range("L40").Value = 2444.4443
$ws.Range("M40").Value = -333964
$ws.Range("N40").Value = -2716.4443
$ws.Range("H122").Value = 2937.375
$ws.Range("I122").Value = 2916.5
$ws.Range("J122").Value = 3000
$ws.Range("K122").Value = 8749.5
$ws.Range("L122").Value = 9000
$ws.Range("M122").Value = -6299.5
$ws.Range("N122").Value = -13900

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H51").Value = 30259.25
$ws.Range("I51").Value = 20000
$ws.Range("J51").Value = 33679
$ws.Range("K51").Value = 20000
$ws.Range("L51").Value = 33679
$ws.Range("M51").Value = -19490
$ws.Range("N51").Value = -34699
$ws.Range("H122").Value = 2362.2354
$ws.Range("I122").Value = 1531.4445
$ws.Range("J122").Value = 3296.875
$ws.Range("K122").Value = 4594.333500000001
$ws.Range("L122").Value = 9890.625
$ws.Range("M122").Value = -2144.333500000001
$ws.Range("N122").Value = -14790.625
